$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ticket's elapsed "on-line" duration (上机时长) is stored as six
# separate numeric components in row 15 (years/months/days/hours/minutes/
# seconds -> B15:G15). Update the hours/minutes/seconds fields:
#   E15 = 上机时长_时 (hours)   : 0  -> 3
#   F15 = 上机时长_分 (minutes) : 22 -> 9
#   G15 = 上机时长_秒 (seconds) : 59 -> 0
# Everything else on the sheet (下机时间, 上网费用, 本次消费, 卡上余额, …)
# is a formula driven off these cells, so it recalculates on its own.
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 0
